$d = $word.ActiveDocument

# Grab the whole document's canonical OOXML once; each top-level <w:p ...>...</w:p>
# block corresponds 1:1 (in order) with $d.Paragraphs.
$full = $d.Content.WordOpenXML
$matches = [regex]::Matches($full, "<w:p\b[^>]*>.*?</w:p>")

$count = $d.Paragraphs.Count
Write-Host ("Paragraphs=" + $count + " XmlBlocks=" + $matches.Count)

for ($i = 0; $i -lt $count; $i++) {
    $paraXml = $matches[$i].Value

    if ($paraXml -notmatch "<w:contextualSpacing\b[^/]*/>") {
        continue
    }

    $newXml = [regex]::Replace($paraXml, "<w:contextualSpacing\b[^/]*/>", "")

    $p = $d.Paragraphs($i + 1)
    $p.Range.InsertXML($newXml)
}

Write-Host "done"
